# Apply translations.xlsx update: add plural/"(s)" variants for param_type, param_group,
# parameter, project and network generic labels (rows 47-60 of final sheet).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new row(s) right after the existing row 50 entry
$ws.Rows.Item(51).Insert()
$ws.Rows.Item(51).Insert()
$ws.Range("A51").Value2 = "networks"
$ws.Range("B51").Value2 = "direct translation"
$ws.Range("C51").Value2 = "Networks"
$ws.Range("D51").Value2 = "Réseaux"
$ws.Range("A52").Value2 = "network(s)"
$ws.Range("B52").Value2 = "direct translation"
$ws.Range("C52").Value2 = "Network(s)"
$ws.Range("D52").Value2 = "Réseau(x)"

# Insert 2 new row(s) right after the existing row 49 entry
$ws.Rows.Item(50).Insert()
$ws.Rows.Item(50).Insert()
$ws.Range("A50").Value2 = "projects"
$ws.Range("B50").Value2 = "direct translation"
$ws.Range("C50").Value2 = "Projects"
$ws.Range("D50").Value2 = "Projets"
$ws.Range("A51").Value2 = "project(s)"
$ws.Range("B51").Value2 = "direct translation"
$ws.Range("C51").Value2 = "Project(s)"
$ws.Range("D51").Value2 = "Projet(s)"

# Insert 2 new row(s) right after the existing row 48 entry
$ws.Rows.Item(49).Insert()
$ws.Rows.Item(49).Insert()
$ws.Range("A49").Value2 = "parameters"
$ws.Range("B49").Value2 = "direct translation"
$ws.Range("C49").Value2 = "Parameters"
$ws.Range("D49").Value2 = "Paramètres"
$ws.Range("A50").Value2 = "parameter(s)"
$ws.Range("B50").Value2 = "direct translation"
$ws.Range("C50").Value2 = "Parameter(s)"
$ws.Range("D50").Value2 = "Paramètre(s)"

# Insert 2 new row(s) right after the existing row 47 entry
$ws.Rows.Item(48).Insert()
$ws.Rows.Item(48).Insert()
$ws.Range("A48").Value2 = "param_groups"
$ws.Range("B48").Value2 = "direct translation"
$ws.Range("C48").Value2 = "Parameter Groups"
$ws.Range("D48").Value2 = "Groupes de paramètres"
$ws.Range("A49").Value2 = "param_group(s)"
$ws.Range("B49").Value2 = "direct translation"
$ws.Range("C49").Value2 = "Parameter Group(s)"
$ws.Range("D49").Value2 = "Groupe(s) de paramètre(s)"

# Insert 2 new row(s) right after the existing row 46 entry
$ws.Rows.Item(47).Insert()
$ws.Rows.Item(47).Insert()
$ws.Range("A47").Value2 = "param_types"
$ws.Range("B47").Value2 = "direct translation"
$ws.Range("C47").Value2 = "Parameter Types"
$ws.Range("D47").Value2 = "Types de paramètres"
$ws.Range("A48").Value2 = "param_type(s)"
$ws.Range("B48").Value2 = "direct translation"
$ws.Range("C48").Value2 = "Parameter Type(s)"
$ws.Range("D48").Value2 = "Type(s) de paramètre(s)"

# Restore view state to match the final edited workbook (scrolled to the new rows
# with D61 as the active cell).
$win = $excel.ActiveWindow
$win.ScrollColumn = 1
$win.ScrollRow = 47
$ws.Range("D61").Select()

